$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column C (rows 2 through 66) from 45202 to 45203 (Förändrad date bumped by one day)
for ($r = 2; $r -le 66; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45202) {
        $cell.Value = 45203
    }
}
